# Free Basics Launch Dates.xlsx - trim data set to the 11 Southeast Asian
# countries now that the mobile-network data collection is complete.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Categorizations": keep only the 11 Southeast Asian countries.
# Drop Bangladesh, Iraq, Maldives, Mongolia, Pakistan, Vanuatu.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Categorizations")
$ws1Rows = @(17, 12, 10, 9, 6, 2)
foreach ($r in $ws1Rows) {
    $ws1.Rows.Item($r).Delete()
}

# ---------------------------------------------------------------------
# Sheet "Asia Pacific Launch Dates": keep only the 5 SE Asia countries
# whose mobile network data is fully populated.
# Drop Bangladesh, Maldives, Mongolia, Pakistan, Vanuatu.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Asia Pacific Launch Dates")
$ws2Rows = @(11, 7, 6, 5, 2)
foreach ($r in $ws2Rows) {
    $ws2.Rows.Item($r).Delete()
}

# ---------------------------------------------------------------------
# Sheet "Metadata": drop the source rows for the removed countries.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Metadata")
$ws3Rows = @(8, 7, 6, 2)
foreach ($r in $ws3Rows) {
    $ws3.Rows.Item($r).Delete()
}

# ---------------------------------------------------------------------
# Re-apply AutoFilter over the now-smaller Categorizations table and
# register the hidden _FilterDatabase defined name that Excel writes
# alongside it.
# ---------------------------------------------------------------------
$ws1.Range("A1:C12").AutoFilter() | Out-Null
$filterName = $ws1.Names.Add("_xlnm._FilterDatabase", "=Categorizations!`$A`$1:`$C`$12")
$filterName.Visible = $false

# ---------------------------------------------------------------------
# Restore each sheet's remembered selection. Order matters: the last
# sheet selected becomes the active tab, and "Metadata" should stay
# active, so it is selected last.
# ---------------------------------------------------------------------
$ws1.Range("A11").Select() | Out-Null
$ws2.Range("C5").Select() | Out-Null
$ws3.Activate() | Out-Null
$ws3.Range("B23").Select() | Out-Null
